$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update task table data ---------------------------------------------
# Row 2 (ID 1, Project A): task changes from "Design UI" to "Testing",
# effort drops to 15, gains a TaskDependency of 2, loses its Progress value.
$ws.Range("D2").Value = "Testing"
$ws.Range("E2").Value = 15
$ws.Range("F2").Value = 2
$ws.Range("H2").ClearContents()

# Row 3: ID renumbered 4 -> 2, project switches from Project B to Project A,
# task becomes "Design UI" with effort 100, loses ProjectDependency, gains Progress 50.
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = "Project A"
$ws.Range("D3").Value = "Design UI"
$ws.Range("E3").Value = 100
$ws.Range("G3").ClearContents()
$ws.Range("H3").Value = 50

# Row 4: ID renumbered 2 -> 3 (task "Implement Backend" stays the same).
$ws.Range("A4").Value = 3

# Row 5: ID renumbered 5 -> 4, TaskDependency updated from 4 to 3.
$ws.Range("A5").Value = 4
$ws.Range("F5").Value = 3

# Row 6: ID renumbered 3 -> 5, project switches from Project A to Project B,
# task becomes "Database Setup" with effort 60, loses TaskDependency, gains ProjectDependency 1.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "Project B"
$ws.Range("D6").Value = "Database Setup"
$ws.Range("E6").Value = 60
$ws.Range("F6").ClearContents()
$ws.Range("G6").Value = 1

# --- Re-sort the table by ID (ascending) instead of EstimatedEffortHours (descending) ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add($ws.Range("A1:A6"))
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# --- Scroll the sheet view so column D is the left-most visible column ---
$ws.Application.ActiveWindow.ScrollColumn = 4
